$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 127, shifting existing rows 127:187 down to 128:188
$ws.Rows(127).Insert()

# Populate the newly inserted row 127 with the new weekly record
$ws.Cells.Item(127, 1).Value = 5
$ws.Cells.Item(127, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(127, 3).Value = "Maule"
$ws.Cells.Item(127, 4).Value = 44466
$ws.Cells.Item(127, 5).Value = 7
$ws.Cells.Item(127, 6).Value = 100114013
$ws.Cells.Item(127, 7).Value = "Zanahoria"
$ws.Cells.Item(127, 8).Value = "Sin especificar"
$ws.Cells.Item(127, 9).Value = "Primera"
$ws.Cells.Item(127, 10).Value = 400
$ws.Cells.Item(127, 11).Value = 6000
$ws.Cells.Item(127, 12).Value = 6000
$ws.Cells.Item(127, 13).Value = 6000
$ws.Cells.Item(127, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(127, 15).Value = "Región de Ñuble"
$ws.Cells.Item(127, 16).Value = 300
$ws.Cells.Item(127, 17).Value = 20
$ws.Cells.Item(127, 18).Value = "Hortaliza"
